# Update the "想去人数" (interested-attendee count) figures in column F
# across the workbook's sheets to match the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 5540
$ws1.Range("F6").Value  = 75
$ws1.Range("F8").Value  = 909
$ws1.Range("F9").Value  = 144
$ws1.Range("F10").Value = 2474
$ws1.Range("F11").Value = 82
$ws1.Range("F12").Value = 101
$ws1.Range("F13").Value = 2
$ws1.Range("F14").Value = 69
$ws1.Range("F15").Value = 5
$ws1.Range("F16").Value = 2317
$ws1.Range("F17").Value = 258

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 100

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 5540
$ws4.Range("F6").Value  = 100
$ws4.Range("F7").Value  = 75
$ws4.Range("F10").Value = 909
$ws4.Range("F11").Value = 144
$ws4.Range("F12").Value = 2474
$ws4.Range("F13").Value = 82
$ws4.Range("F14").Value = 101
$ws4.Range("F15").Value = 2
$ws4.Range("F17").Value = 69
$ws4.Range("F18").Value = 5
$ws4.Range("F19").Value = 2317
$ws4.Range("F20").Value = 258
